$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New rows 7 & 8: more vertical 5ml runs, same day (filenames only for now) ---
$ws.Range("A7").Value = "D20151103T161152"
$ws.Range("A8").Value = "D20151103T163411"

# --- Note that one of the earlier vertical runs used all signals (not just a subset) ---
$ws.Range("F3").Value = "use all signals"

# --- New columns for the IFCB run-stats captured alongside CellConc ---
$ws.Range("E1:H1").EntireColumn.Insert()
$ws.Range("E1").Value = "ml_analyzed"
$ws.Range("F1").Value = "runtime"
$ws.Range("G1").Value = "inhibittime"
$ws.Range("H1").Value = "numtriggers"

# Fill in the measurement data that was collected for row 6 (previously only had a comment)
$ws.Range("D6").Value = 637.70000000000005
$ws.Range("E6").Value = 4.0865
$ws.Range("F6").Value = 1198.04
$ws.Range("G6").Value = 202.27199999999999
$ws.Range("H6").Value = 2606

# --- Finish filling in rows 7 & 8 ---
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = "V"
$ws.Range("D7").Value = 636.6
$ws.Range("E7").Value = 4.0888
$ws.Range("F7").Value = 1198
$ws.Range("G7").Value = 204.11
$ws.Range("H7").Value = 2603
$ws.Range("I7").Value = "9um beads"

$ws.Range("B8").Value = 5
$ws.Range("C8").Value = "V"
$ws.Range("D8").Value = 635.6
$ws.Range("E8").Value = 4.0873999999999997
$ws.Range("F8").Value = 1198
$ws.Range("G8").Value = 201.6
$ws.Range("H8").Value = 2598
$ws.Range("I8").Value = "9um beads"

# --- Row 9: concentration run cut short, about to switch to horz 5ml runs ---
$ws.Range("A9").Value = "D20151103T165631"
$ws.Range("B9").Value = "<5"
$ws.Range("C9").Value = "V"
$ws.Range("D9").Value = 600
$ws.Range("E9").Value = 1.1843999999999999
$ws.Range("F9").Value = 351.2
$ws.Range("G9").Value = 56.13
$ws.Range("H9").Value = 711
$ws.Range("I9").Value = "sampled ended early because concentration consistent enough to move on to next stage: horz 5ml runs"
$ws.Range("J9").Value = "9um beads, use all signals"

# --- Mark the rest of the "9um beads, use all signals" rows accordingly ---
$ws.Range("K2").Value = "use all signals"
$ws.Range("K4").Value = "use all signals"
$ws.Range("J5").Value = "use all signals"
$ws.Range("J6").Value = "use all signals"
$ws.Range("J7").Value = "use all signals"
$ws.Range("J8").Value = "use all signals"

# --- Row 10: first file of the next (horizontal) stage ---
$ws.Range("B10").Value = 5
$ws.Range("C10").Value = "H"

# Column widths for the four new columns, matching the other data columns
$ws.Range("E1:H1").ColumnWidth = 12.14

# Update selection to reflect where the author left off
$ws.Range("A10").Select()
